# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet with newly scraped figures. Column D values are plain numeric-
# looking text (e.g. "27.988.12", "0.5128") that must stay stored as
# text, so we force a text number format before assigning each one and
# then clear the format again so the cell's style matches the rest of
# the sheet (unstyled) once the text is safely in place. Column E
# values ("  -0.27%  ") are never auto-converted by Excel, so they are
# assigned directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = "27.988.12"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = "  -0.27%  "
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = "1.857.41"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = "  -0.86%  "
$ws.Range('E4').Value = "  +0.09%  "
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = "312.60"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = "  -0.30%  "
$ws.Range('E6').Value = "  +0.07%  "
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = "0.5128"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = "  +1.30%  "
$ws.Range('E8').Value = "  -0.37%  "
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = "0.08237"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = "  -8.24%  "
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = "1.109"
$ws.Range('D10').ClearFormats()
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = "41.51"
$ws.Range('D11').ClearFormats()
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = "6.183"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = "  -2.41%  "
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = "20.55"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = "  -0.75%  "
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = "1.863.57"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = "  -0.19%  "
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = "7.239"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = "  +0.65%  "
$ws.Range('E16').Value = "  +0.04%  "
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = "0.00001097"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = "  -1.03%  "
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = "90.52"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = "  -0.55%  "
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = "0.06649"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = "  +0.85%  "
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = "17.68"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = "  -2.55%  "
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = "5.997"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = "  -1.82%  "
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = "28.010.97"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = "  -0.25%  "
$ws.Range('E24').Value = "  -3.16%  "
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = "2.246"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = "  -1.61%  "
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = "2.073.16"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = "  -0.65%  "
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = "2.509"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = "  -1.16%  "
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = "157.97"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = "  +0.34%  "
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = "20.43"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = "  -1.59%  "
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = "124.48"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = "  -1.59%  "
$ws.Range('E31').Value = "  +1.10%  "
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = "1.029"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = "  -2.83%  "
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = "5.957"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = "  +6.15%  "
$ws.Range('E34').Value = "  -0.12%  "
$ws.Range('E35').Value = "  -3.10%  "
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = "0.02412"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = "  -0.56%  "
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = "0.06487"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = "  -1.32%  "
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = "0.2168"
$ws.Range('D38').ClearFormats()
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = "0.6520"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = "  +2.15%  "
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = "1.194"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = "  -1.04%  "
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = "5.016"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = "  +2.20%  "
$ws.Range('E42').Value = "  -3.49%  "
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = "11.14"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = "  -2.54%  "
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = "0.6141"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = "  +2.20%  "
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = "12.99"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = "  -1.24%  "
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = "1.279"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = "  +0.26%  "
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = "3.659"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = "  -0.32%  "
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = "1.216"
$ws.Range('D49').ClearFormats()
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = "120.29"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = "  -0.85%  "
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = "78.31"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = "  -1.68%  "
